$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A4").Value = -21.774
$ws.Range("E4").Value = 13.418

$ws.Range("E5").Value = 13.429

$ws.Range("A6").Value = -21.059
$ws.Range("E6").Value = 13.363

$ws.Range("A7").Value = -21.089

$ws.Range("A8").Value = -20.846
$ws.Range("E8").Value = 13.486

$ws.Range("A16").Value = -20.65
$ws.Range("E16").Value = 13.068

$ws.Range("A20").Value = -22.138

$ws.Range("A21").Value = -21.14

$ws.Range("E22").Value = 13.342
